# Regenerate merged AHB files
# 1) Rename the header labels in row 1 from the *_old/*_new suffix scheme
#    to the *_FV2310/*_FV2404 scheme.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headerMap = @{
    "A1" = "Segmentname_FV2310"
    "B1" = "Segmentgruppe_FV2310"
    "C1" = "Segment_FV2310"
    "D1" = "Datenelement_FV2310"
    "E1" = "Segment ID_FV2310"
    "F1" = "Code_FV2310"
    "G1" = "Qualifier_FV2310"
    "H1" = "Beschreibung_FV2310"
    "I1" = "Bedingungsausdruck_FV2310"
    "J1" = "Bedingung_FV2310"
    "K1" = "diff"
    "L1" = "Segmentname_FV2404"
    "M1" = "Segmentgruppe_FV2404"
    "N1" = "Segment_FV2404"
    "O1" = "Datenelement_FV2404"
    "P1" = "Segment ID_FV2404"
    "Q1" = "Code_FV2404"
    "R1" = "Qualifier_FV2404"
    "S1" = "Beschreibung_FV2404"
    "T1" = "Bedingungsausdruck_FV2404"
    "U1" = "Bedingung_FV2404"
}

foreach ($addr in $headerMap.Keys) {
    $ws.Range($addr).Value = $headerMap[$addr]
}

# 2) Turn the whole used range into an Excel Table ("Table1") so the sheet
#    gets an autofilter + structured table definition, reusing the header
#    row that is already present.
$rng = $ws.Range("A1:U57")
$tbl = $ws.ListObjects.Add(1, $rng, $false, 1, $null)
$tbl.Name = "Table1"

# 3) Freeze the header row (split/freeze at row 2, i.e. 1 row frozen).
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
